# fix input and output conventions
#
# The sheet previously carried a redundant leading "Topic" column (always
# "Verizon") ahead of the Comment/Result/Sentiment/Confidence columns. That
# column is removed entirely, shifting Comment into column A, Result into
# column B, and Sentiment/Confidence into columns C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "Topic" column (A) - everything shifts one column left.
$ws.Columns("A").Delete()

# Match the saved selection/active cell from the edit.
$ws.Range("C15").Select()
